$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")

# Row 9 - "Realizar descripción de CU 01 y CU 02": log 1 estimated hour,
# with 1 hour consumed on Día 6 (column W)
$ws.Range("G9").Value = 1
$ws.Range("W9").Value = 1

# Row 11 - "Realizar mockup de ventana principal de profesor": 1 hour
# consumed on Día 6 (column W); I11 (Día 1 remaining) is overwritten
# directly with its resulting value
$ws.Range("I11").Value = 1
$ws.Range("W11").Value = 1

# Re-merge the day-header cells at the tail end of the header row; this
# matches the merge order left behind after the edit session
$mergedHeaderRanges = @("AZ4:BA4", "AO4:AP4", "AR4:AS4", "AU4:AV4", "AX4:AY4")
foreach ($r in $mergedHeaderRanges) {
    $ws.Range($r).UnMerge()
}
foreach ($r in $mergedHeaderRanges) {
    $ws.Range($r).Merge()
}

# Reflect the final UI state left behind by the edit session: zoomed to
# 70% in Page Break Preview, with the last active selection on W12
$win = $excel.ActiveWindow
$win.Zoom = 70
$ws.Range("W12").Select()
